# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" period codes (shared strings 1709..1902) are put in
# ascending order (previously descending, 1902..1709), and the "Valor Mora"
# (F) / "Salario Basico" (G) figures are refreshed for every worker row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("1709","1710","1711","1712","1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902")
$valorMora     = @(29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,31249,31249,31249,31249,31249,31249)
$salarioBasico = @(828116,828116,828116,828116,828116,828116,828116,828116,828116,828116,828116,828116,828116,828116,828116,828116,828116,828116)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valorMora[$i]
    $ws.Range("G$row").Value = $salarioBasico[$i]
}
